$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.322699999999997
$ws.Range("B3").Value = 5.502299999999992
$ws.Range("B5").Value = 4.798100000000003
$ws.Range("C7").Value = -11.80969999999999
$ws.Range("A9").Value = -20.40919999999998
$ws.Range("C9").Value = -11.93670000000001
$ws.Range("B11").Value = 5.340899999999996
$ws.Range("B12").Value = 5.670899999999997
$ws.Range("A13").Value = -21.81870000000002
$ws.Range("A16").Value = -20.1366
$ws.Range("A18").Value = -22.7221
$ws.Range("A20").Value = -22.03500000000002
$ws.Range("B21").Value = 5.368999999999995
$ws.Range("C21").Value = -13.52180000000001
